$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Actual Result (E3) gets a new session URL ---
$ws.Range("E3").Value = "http://localhost:8080/home/index;jsessionid=647D78420C01EDFF052B79870B73076D"

# --- Row 4: Actual Result (E4) becomes a long multi-line Selenium error message,
#     and Status (F4) flips from PASS to FAIL ---
$e4Text = @"
Lỗi: no such window: target window already closed
from unknown error: web view not found
  (Session info: chrome=145.0.7632.76)
Build info: version: '4.14.1', revision: '03f8ede370'
System info: os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '17.0.10'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [b5d05d6c08870a977d15179599401354, sendKeysToElement {id=f.F24A41F9119A8E04E81E2562BC6393EE.d.48E05C2B15138ED5E3632BC7F8322B81.e.4, value=[Ljava.lang.CharSequence;@68868328}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.76, chrome: {chromedriverVersion: 145.0.7632.77 (da516187054a..., userDataDir: C:\Users\DELL\AppData\Local...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:50697}, goog:processID: 1616, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:50697/devtoo..., se:cdpVersion: 145.0.7632.76, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Element: [[ChromeDriver: chrome on windows (b5d05d6c08870a977d15179599401354)] -> name: password]
Session ID: b5d05d6c08870a977d15179599401354
"@

# Setting a multi-line value directly on E4 triggers an automatic row-height
# autofit side effect that the source workbook never had (row 4 keeps its
# original ht="15" with no customHeight flag). Route the value through a
# scratch cell + copy/paste-special so the destination's row metrics are left
# untouched, then clean the scratch cell/row back up.
$scratch = $ws.Cells.Item(20, 20)
$scratch.Value = $e4Text
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()
$ws.Rows.Item(20).EntireRow.Delete()
$excel.CutCopyMode = $false

$ws.Range("F4").Value = "FAIL"
